$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "as20232248@sva.edu.eg"
$ws.Range("B67").Value = "https://hoda3225.github.io/my-first-webpage/"
$ws.Range("A68").Value = "ys20231357@sva.edu.eg"
$ws.Range("B68").Value = "https://goo782.github.io/My-page/"
